$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.105724
$ws.Range("H2").Value = 30.317172
$ws.Range("I2").Value = 0.5504853801993582
$ws.Range("J2").Value = 0.5504853801993582
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.33599166666667
$ws.Range("N2").Value = 52.007975
$ws.Range("O2").Value = 0.4573561888773979
$ws.Range("P2").Value = 0.4573561888773979
$ws.Range("Q2").Value = 175.1927470496334
$ws.Range("R2").Value = 1576.7347234467
$ws.Range("S2").Value = 0.2517678955207038
$ws.Range("T2").Value = 0.2517678955207038
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.105724
$ws.Range("H3").Value = 30.317172
$ws.Range("I3").Value = 0.5504853801993582
$ws.Range("J3").Value = 0.5504853801993582
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.077707333333334
$ws.Range("N3").Value = 27.233122
$ws.Range("O3").Value = 0.2394870573052156
$ws.Range("P3").Value = 0.2394870573052156
$ws.Range("Q3").Value = 91.73680486344269
$ws.Range("R3").Value = 825.631243770984
$ws.Range("S3").Value = 0.1318341237934871
$ws.Range("T3").Value = 0.1318341237934871
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.105724
$ws.Range("H4").Value = 30.317172
$ws.Range("I4").Value = 0.5504853801993582
$ws.Range("J4").Value = 0.5504853801993582
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.491094
$ws.Range("N4").Value = 34.473282
$ws.Range("O4").Value = 0.3031567538173866
$ws.Range("P4").Value = 0.3031567538173866
$ws.Range("Q4").Value = 116.125824422056
$ws.Range("R4").Value = 1045.132419798504
$ws.Range("S4").Value = 0.1668833608851673
$ws.Range("T4").Value = 0.1668833608851673
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.009378000000001
$ws.Range("H5").Value = 15.028134
$ws.Range("I5").Value = 0.2728740021884925
$ws.Range("J5").Value = 0.2728740021884924
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.33599166666667
$ws.Range("N5").Value = 52.007975
$ws.Range("O5").Value = 0.4573561888773979
$ws.Range("P5").Value = 0.4573561888773979
$ws.Range("Q5").Value = 86.84253526318335
$ws.Range("R5").Value = 781.5828173686501
$ws.Range("S5").Value = 0.1248006136846517
$ws.Range("T5").Value = 0.1248006136846516
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.009378000000001
$ws.Range("H6").Value = 15.028134
$ws.Range("I6").Value = 0.2728740021884925
$ws.Range("J6").Value = 0.2728740021884924
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.077707333333334
$ws.Range("N6").Value = 27.233122
$ws.Range("O6").Value = 0.2394870573052156
$ws.Range("P6").Value = 0.2394870573052156
$ws.Range("Q6").Value = 45.47366740603868
$ws.Range("R6").Value = 409.2630066543481
$ws.Range("S6").Value = 0.06534979179921903
$ws.Range("T6").Value = 0.06534979179921901
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.009378000000001
$ws.Range("H7").Value = 15.028134
$ws.Range("I7").Value = 0.2728740021884925
$ws.Range("J7").Value = 0.2728740021884924
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.491094
$ws.Range("N7").Value = 34.473282
$ws.Range("O7").Value = 0.3031567538173866
$ws.Range("P7").Value = 0.3031567538173866
$ws.Range("Q7").Value = 57.563233479532
$ws.Range("R7").Value = 518.069101315788
$ws.Range("S7").Value = 0.08272359670462183
$ws.Range("T7").Value = 0.08272359670462182
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.242740666666667
$ws.Range("H8").Value = 9.728222000000001
$ws.Range("I8").Value = 0.1766406176121494
$ws.Range("J8").Value = 0.1766406176121493
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.33599166666667
$ws.Range("N8").Value = 52.007975
$ws.Range("O8").Value = 0.4573561888773979
$ws.Range("P8").Value = 0.4573561888773979
$ws.Range("Q8").Value = 56.21612517449445
$ws.Range("R8").Value = 505.94512657045
$ws.Range("S8").Value = 0.0807876796720424
$ws.Range("T8").Value = 0.08078767967204238
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.242740666666667
$ws.Range("H9").Value = 9.728222000000001
$ws.Range("I9").Value = 0.1766406176121494
$ws.Range("J9").Value = 0.1766406176121493
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.077707333333334
$ws.Range("N9").Value = 27.233122
$ws.Range("O9").Value = 0.2394870573052156
$ws.Range("P9").Value = 0.2394870573052156
$ws.Range("Q9").Value = 29.43665072989823
$ws.Range("R9").Value = 264.9298565690841
$ws.Range("S9").Value = 0.04230314171250949
$ws.Range("T9").Value = 0.04230314171250948
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.242740666666667
$ws.Range("H10").Value = 9.728222000000001
$ws.Range("I10").Value = 0.1766406176121494
$ws.Range("J10").Value = 0.1766406176121493
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 11.491094
$ws.Range("N10").Value = 34.473282
$ws.Range("O10").Value = 0.3031567538173866
$ws.Range("P10").Value = 0.3031567538173866
$ws.Range("Q10").Value = 37.26263781828933
$ws.Range("R10").Value = 335.363740364604
$ws.Range("S10").Value = 0.05354979622759749
$ws.Range("T10").Value = 0.05354979622759748

